# LOM3108.xlsx: add a new "Docentes responsáveis:" (responsible faculty)
# section right after the existing "Objectives:" row.
#
# Before:  row 11 = "Objectives:" ... row 12 = "Programa resumido:" (+2 data cols)
# After:   row 11 = "Objectives:"
#          row 12 = "Docentes responsáveis:"           (label only, col A)
#          row 13 = "7459752 - Maria Ismenia Sodero Toledo Faria"   (cols B/C)
#          row 14 = "2166002 - Sandra Giacomin Schneider"           (cols B/C)
#          row 15 = "1922320 - Sebastiao Ribeiro"                   (cols B/C)
#          row 16 = "Programa resumido:" ... (everything below shifts down by 4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 12..15, pushing the old row 12 ("Programa resumido:")
# and everything after it down to row 16 onward (old row 12 kept its original
# ht="60" customHeight formatting, now living on row 16).
$ws.Rows("12:15").Insert()

# The insert copies row 11's formatting down into the new rows, which leaves a
# stray formatted-but-empty cell in column A for rows 13-15 (no A-column entry
# belongs on those rows). Clear them out completely.
$ws.Range("A13:A15").Clear()

# New section label (bold, column A only - same style as the other section
# headers such as "Objectives:"/"Programa resumido:").
$ws.Range("A12").Value = "Docentes responsáveis:"

# Three faculty rows, duplicated into both the "current" (B) and "modified"
# (C) syllabus columns, same as every other two-column data row in the sheet.
$ws.Range("B13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("B14").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C14").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"

# Give the new B/C cells the same formatting (wrap text, top-aligned; C in
# red font) as the rest of the table by copying it from row 16, which still
# carries the original column styles.
$ws.Range("B16:C16").Copy()
$ws.Range("B13:C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
